$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34: turn into the last row of its text-block group (bottom border) ---
# Copy formatting (incl. border) from the existing last-row pattern (row 30) onto row 34.
# Cell values (B34 number, C34/D34/E34 existing shared strings) are unchanged.
$ws.Range("A30:E30").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)

# --- Row 35 (new): SCRIPT/T01P01A/um2205.ssb block, first row ---
$ws.Range("A35").Value = 'SCRIPT/T01P01A/um2205.ssb'
$ws.Range("B35").Value = 416
$ws.Range("C35").Value = ' Hey, hey! [CS:P]Brine Cave[CR] has some\nnasty enemies on the loose.'
$ws.Range("D35").Value = ' Эй, эй! В [CS:P]Пещере у Моря[CR] нас\nподжидают очень опасные враги.'
$ws.Range("E35").Value = ' Üê, üê! Â [CS:P]Ðåþåñå ô Íïñÿ[CR] îàò\nðïäçéäàýó ïœåîû ïðàòîúå âñàãé.'
$ws.Rows.Item(35).RowHeight = 43.2

# --- Row 36 (new): same block, last row (bottom border) ---
$ws.Range("A30:E30").Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)
$ws.Range("B36").Value = 419
$ws.Range("C36").Value = ' Everyone\''d better get stocked up\nwith proper equipment, hey, hey.'
$ws.Range("D36").Value = ' Нам всем нужно запастись\nподходящим снаряжением, эй, эй.'
$ws.Range("E36").Value = ' Îàí âòåí îôçîï èàðàòóéòû\nðïäöïäÿþéí òîàñÿçåîéåí, üê, üê.'
$ws.Rows.Item(36).RowHeight = 21.6

# --- Row 37 (new): SCRIPT/T01P02A/um2406.ssb block, first row ---
$ws.Range("A37").Value = 'SCRIPT/T01P02A/um2406.ssb'
$ws.Range("B37").Value = 394
$ws.Range("C37").Value = ' Hey, hey! I believe in you!'
$ws.Range("D37").Value = ' Эй, эй! Я верю в вас!'
$ws.Range("E37").Value = ' Üê, üê! Ÿ âåñý â âàò!'
$ws.Rows.Item(37).RowHeight = 43.2

# --- Row 38 (new): same block, last row (regular style, no bottom border, no A cell) ---
$ws.Range("B38").Value = 397
$ws.Range("C38").Value = ' I know you\''ll get things done,\nhey, hey!'
$ws.Range("D38").Value = ' Я знаю, что у вас всё получится,\nэй, эй!'
$ws.Range("E38").Value = ' Ÿ èîàý, œóï ô âàò âòæ ðïìôœéóòÿ,\nüê, üê!'
$ws.Rows.Item(38).RowHeight = 21.6

# --- Update the saved view: active cell / selection moves to E38 ---
$ws.Range("E38").Select()
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
